# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# columns for the 4821b8ad... file row on both the zh-cn and de-de handback
# status sheets, reflecting the newly generated handback report timestamps.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 3 is the 4821b8ad...zh-cn.xlf entry
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-21 16:49:29"
$wsZhCn.Range("H3").Value = "2016-03-21 16:49:50"

# de-de sheet: row 3 is the 4821b8ad...de-de.xlf entry
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-21 16:49:33"
$wsDeDe.Range("H3").Value = "2016-03-21 16:49:56"
